$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) from 45172 to 45175 for rows 2..477
for ($r = 2; $r -le 477; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45172) {
        $cell.Value = 45175
    }
}

# Add row 478
$ws.Cells.Item(478, 1).Value = "A 40894-2023"
$ws.Cells.Item(478, 2).Value = 45173
$ws.Cells.Item(478, 3).Value = 45175
$ws.Cells.Item(478, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(478, 5).Value = "EKSJÖ"
$ws.Cells.Item(478, 7).Value = 1.5
$ws.Cells.Item(478, 8).Value = 0
$ws.Cells.Item(478, 9).Value = 0
$ws.Cells.Item(478, 10).Value = 0
$ws.Cells.Item(478, 11).Value = 0
$ws.Cells.Item(478, 12).Value = 0
$ws.Cells.Item(478, 13).Value = 0
$ws.Cells.Item(478, 14).Value = 0
$ws.Cells.Item(478, 15).Value = 0
$ws.Cells.Item(478, 16).Value = 0
$ws.Cells.Item(478, 17).Value = 0

# Add row 479
$ws.Cells.Item(479, 1).Value = "A 41157-2023"
$ws.Cells.Item(479, 2).Value = 45174
$ws.Cells.Item(479, 3).Value = 45175
$ws.Cells.Item(479, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(479, 5).Value = "EKSJÖ"
$ws.Cells.Item(479, 7).Value = 1.9
$ws.Cells.Item(479, 8).Value = 0
$ws.Cells.Item(479, 9).Value = 0
$ws.Cells.Item(479, 10).Value = 0
$ws.Cells.Item(479, 11).Value = 0
$ws.Cells.Item(479, 12).Value = 0
$ws.Cells.Item(479, 13).Value = 0
$ws.Cells.Item(479, 14).Value = 0
$ws.Cells.Item(479, 15).Value = 0
$ws.Cells.Item(479, 16).Value = 0
$ws.Cells.Item(479, 17).Value = 0

# Copy formatting (style) for B, C and R columns from row 477
$ws.Range("B477:C477").Copy()
$ws.Range("B478:C479").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R477").Copy()
$ws.Range("R478:R479").PasteSpecial(-4122)

# Row heights: row 477 and 478 get explicit 15pt custom height, 479 doesn't
$ws.Rows.Item(477).RowHeight = 15
$ws.Rows.Item(478).RowHeight = 15
